$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 2500
$ws.Range("I34").Value = 2500
$ws.Range("K34").Value = 2500
$ws.Range("M34").Value = -2297
$ws.Range("H36").Value = 2500
$ws.Range("I36").Value = 2500
$ws.Range("K36").Value = 2500
$ws.Range("M36").Value = -1785
$ws.Range("H76").Value = 3351999.8
$ws.Range("I76").Value = 5859601
$ws.Range("K76").Value = 5859601
$ws.Range("M76").Value = -5859286
$ws.Range("H79").Value = 3351999.8
$ws.Range("I79").Value = 5859601
$ws.Range("K79").Value = 5859601
$ws.Range("M79").Value = -5858509
$ws.Range("H86").Value = 1666.6666
$ws.Range("I86").Value = 1500
$ws.Range("K86").Value = 1500
$ws.Range("M86").Value = -377
$ws.Range("H89").Value = 1666.6666
$ws.Range("I89").Value = 1500
$ws.Range("K89").Value = 7500
$ws.Range("M89").Value = -1884
$ws.Range("H129").Value = 855.3409
$ws.Range("I129").Value = 693.1667
$ws.Range("J129").Value = 880.9474
$ws.Range("K129").Value = 2079.5001
$ws.Range("L129").Value = 2642.8422
$ws.Range("M129").Value = 2920.4999
$ws.Range("N129").Value = -12642.8422
$ws.Range("H138").Value = 2856.84
$ws.Range("I138").Value = 2512.3333
$ws.Range("J138").Value = 3373.6
$ws.Range("K138").Value = 7536.999899999999
$ws.Range("L138").Value = 10120.8
$ws.Range("M138").Value = -2396.999899999999
$ws.Range("N138").Value = -20400.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6686.2
$ws.Range("I61").Value = 3799.5
$ws.Range("J61").Value = 7407.875
$ws.Range("K61").Value = 3799.5
$ws.Range("L61").Value = 7407.875
$ws.Range("M61").Value = -3587.5
$ws.Range("N61").Value = -7831.875
$ws.Range("H88").Value = 4487.5
$ws.Range("I88").Value = 2466.6667
$ws.Range("K88").Value = 2466.6667
$ws.Range("M88").Value = -2060.6667
$ws.Range("H91").Value = 4487.5
$ws.Range("I91").Value = 2466.6667
$ws.Range("K91").Value = 2466.6667
$ws.Range("M91").Value = -1062.6667
$ws.Range("H122").Value = 1836.5
$ws.Range("I122").Value = 1836.5
$ws.Range("K122").Value = 5509.5
$ws.Range("M122").Value = -3059.5
$ws.Range("H132").Value = 2233.4583
$ws.Range("I132").Value = 1817.9166
$ws.Range("K132").Value = 5453.7498
$ws.Range("M132").Value = -2923.7498
$ws.Range("H136").Value = 6686.2
$ws.Range("I136").Value = 3799.5
$ws.Range("J136").Value = 7407.875
$ws.Range("K136").Value = 11398.5
$ws.Range("L136").Value = 22223.625
$ws.Range("M136").Value = -8848.5
$ws.Range("N136").Value = -27323.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 136300.67
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 184773.64
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 184773.64
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -187019.64
$ws.Range("H89").Value = 136300.67
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 184773.64
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 923868.2000000001
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -935100.2000000001
$ws.Range("H105").Value = 2005.5862
$ws.Range("I105").Value = 1948.5769
$ws.Range("K105").Value = 1948.5769
$ws.Range("M105").Value = -201.5769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9552.632
$ws.Range("J4").Value = 9552.632
$ws.Range("L4").Value = 9552.632
$ws.Range("N4").Value = -9776.632
$ws.Range("H16").Value = 847.5
$ws.Range("I16").Value = 809.75
$ws.Range("K16").Value = 809.75
$ws.Range("M16").Value = -522.75
$ws.Range("H58").Value = 2417948.8
$ws.Range("I58").Value = 4349695.5
$ws.Range("K58").Value = 4349695.5
$ws.Range("M58").Value = -4349492.5
$ws.Range("H113").Value = 847.5
$ws.Range("I113").Value = 809.75
$ws.Range("K113").Value = 809.75
$ws.Range("M113").Value = 1360.25
$ws.Range("H132").Value = 2700.5264
$ws.Range("I132").Value = 1263.5454
$ws.Range("J132").Value = 4676.375
$ws.Range("K132").Value = 3790.6362
$ws.Range("L132").Value = 14029.125
$ws.Range("M132").Value = -1260.6362
$ws.Range("N132").Value = -19089.125
$ws.Range("H136").Value = 2417948.8
$ws.Range("I136").Value = 4349695.5
$ws.Range("K136").Value = 13049086.5
$ws.Range("M136").Value = -13046536.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 889.73334
$ws.Range("I107").Value = 500
$ws.Range("J107").Value = 917.5714
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 2752.7142
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -6592.7142
$ws.Range("H122").Value = 992.8182
$ws.Range("J122").Value = 1225.8572
$ws.Range("L122").Value = 11032.7148
$ws.Range("N122").Value = -15932.7148
$ws.Range("H131").Value = 783.9
$ws.Range("J131").Value = 792.6804
$ws.Range("L131").Value = 2378.0412
$ws.Range("N131").Value = -12458.0412
$ws.Range("H141").Value = 3562
$ws.Range("I141").Value = 3562
$ws.Range("K141").Value = 10686
$ws.Range("M141").Value = -5506

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13000.777
$ws.Range("I70").Value = 23747.5
$ws.Range("J70").Value = 4403.4
$ws.Range("K70").Value = 23747.5
$ws.Range("L70").Value = 4403.4
$ws.Range("M70").Value = -23477.5
$ws.Range("N70").Value = -4943.4
$ws.Range("H73").Value = 13000.777
$ws.Range("I73").Value = 23747.5
$ws.Range("J73").Value = 4403.4
$ws.Range("K73").Value = 23747.5
$ws.Range("L73").Value = 4403.4
$ws.Range("M73").Value = -22811.5
$ws.Range("N73").Value = -6275.4
$ws.Range("H102").Value = 2632.4
$ws.Range("I102").Value = 2676.5334
$ws.Range("K102").Value = 2676.5334
$ws.Range("M102").Value = -1054.5334
$ws.Range("H132").Value = 1833824.6
$ws.Range("I132").Value = 3206803.2
$ws.Range("J132").Value = 3186.4443
$ws.Range("K132").Value = 9620409.600000001
$ws.Range("L132").Value = 9559.332900000001
$ws.Range("M132").Value = -9617879.600000001
$ws.Range("N132").Value = -14619.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 7378.4287
$ws.Range("I32").Value = 6381.8
$ws.Range("K32").Value = 6381.8
$ws.Range("M32").Value = -6064.8
$ws.Range("H46").Value = 2316.3333
$ws.Range("I46").Value = 1411.5
$ws.Range("J46").Value = 2919.5557
$ws.Range("K46").Value = 1411.5
$ws.Range("L46").Value = 2919.5557
$ws.Range("M46").Value = -1223.5
$ws.Range("N46").Value = -3295.5557
$ws.Range("H136").Value = 2345.2144
$ws.Range("I136").Value = 1388.8572
$ws.Range("K136").Value = 4166.571599999999
$ws.Range("M136").Value = -1616.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1405.579
$ws.Range("I132").Value = 1159.2354
$ws.Range("K132").Value = 3477.7062
$ws.Range("M132").Value = -947.7062000000001
$ws.Range("H136").Value = 25255510
$ws.Range("I136").Value = 34725300
$ws.Range("J136").Value = 2733.3333
$ws.Range("K136").Value = 104175900
$ws.Range("L136").Value = 8199.999899999999
$ws.Range("M136").Value = -104173350
$ws.Range("N136").Value = -13299.9999
